$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# Column G (status) width: XML stored width 18 -> 10.
# COM ColumnWidth <-> stored XML width differs by a fixed padding offset
# (observed offset = 5/6 char), so request 10 - 5/6 to land exactly on 10.
$ws.Columns.Item(7).ColumnWidth = 10 - 5/6

# Games that were in-progress ("13:36 - 2nd Half") have since finished;
# update their status to "Final".
$rows = @(8, 10, 12, 14, 16, 17, 27, 37, 40, 41, 50, 51, 55, 58, 60, 62, 65, 72, 73)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = "Final"
}
